$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" header in H1, matching the style of the existing header cells (e.g. G1 "sum")
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# New data values in column H for the existing data rows
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
